$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the checkmark values out of the two result tables (B4:E7 and
# B12:E13, B15:E15) while keeping their existing cell formatting/style.
$ws.Range("B4:E7").ClearContents()
$ws.Range("B12:E13").ClearContents()
$ws.Range("B15:E15").ClearContents()

# Row 14 ("Fixed file") previously held "Fails" labels (red, centered
# style) in B14:E14 -- remove those cells entirely (contents + formatting)
# so the now-unused "Fails" shared string is dropped as well.
$ws.Range("B14:E14").Clear()

# Update the selection to match the new active range.
$ws.Range("B12:E15").Select()
